$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds a stacked table of Região / Ano / Índice de Gini rows,
# grouped by region (Brasil, Nordeste, Sergipe) each spanning years
# 2012-2023. We are adding a 2024 data point to each region's block,
# which means inserting one row at the end of each region's block and
# shifting everything below it down.

# 1) Insert a new row right before the Nordeste block (currently starting
#    at row 14) to hold the new "Brasil 2024" record.
$ws.Rows("14:14").Insert()

# 2) Insert a new row right before the Sergipe block. After step 1 the
#    Sergipe block starts at row 27, so this row becomes "Nordeste 2024".
$ws.Rows("27:27").Insert()

# 3) Insert a new row right after the (now shifted) Sergipe block, which
#    ends at row 39, to hold "Sergipe 2024".
$ws.Rows("40:40").Insert()

# Fill in the three new rows.
$ws.Range("A14").Value = "Brasil"
$ws.Range("B14").Value = 2024
$ws.Range("C14").Value = 0.506

$ws.Range("A27").Value = "Nordeste"
$ws.Range("B27").Value = 2024
$ws.Range("C27").Value = 0.502

$ws.Range("A40").Value = "Sergipe"
$ws.Range("B40").Value = 2024
$ws.Range("C40").Value = 0.5
